$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = 'Cluster Name'
$ws.Range("B1").Value = 'Active cases'

# Data rows 2-35
$ws.Range("A2").Value = '3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North'
$ws.Range("B2").Value = 13
$ws.Range("A3").Value = '3376 Royal Freemasons Coppin Centre Melbourne'
$ws.Range("B3").Value = 11
$ws.Range("A4").Value = '3601 Baptcare Westhaven community'
$ws.Range("B4").Value = 28
$ws.Range("A5").Value = '3653 Fronditha Thalpori St Albans Aged Care'
$ws.Range("B5").Value = 39
$ws.Range("A6").Value = '44121 Wallan Primary School Wallan'
$ws.Range("B6").Value = 17
$ws.Range("A7").Value = '44165 Greenvale Primary School'
$ws.Range("B7").Value = 22
$ws.Range("A8").Value = '44234 Lucknow Primary School Bairnsdale'
$ws.Range("B8").Value = 25
$ws.Range("A9").Value = '44321 Maiden Gully Primary School Maiden Gully'
$ws.Range("B9").Value = 12
$ws.Range("A10").Value = '44395 Buln Buln Primary School'
$ws.Range("B10").Value = 12
$ws.Range("A11").Value = '44701 Hampton Park Primary School Hampton Park'
$ws.Range("B11").Value = 10
$ws.Range("A12").Value = '44811 Dandenong North Primary School Dandenong'
$ws.Range("B12").Value = 24
$ws.Range("A13").Value = '44853 St Albans North Primary School'
$ws.Range("B13").Value = 12
$ws.Range("A14").Value = '45158 Rowellyn Park Primary School Carrum Downs'
$ws.Range("B14").Value = 16
$ws.Range("A15").Value = '45249 Creekside K-9 College Caroline Springs'
$ws.Range("B15").Value = 16
$ws.Range("A16").Value = '45695 Sacred Heart Primary School Yarrawonga'
$ws.Range("B16").Value = 30
$ws.Range("A17").Value = '4574 Village Glen Aged Care Residences Mornington'
$ws.Range("B17").Value = 20
$ws.Range("A18").Value = '45809 St Finbar''s Primary School Brighton East'
$ws.Range("B18").Value = 11
$ws.Range("A19").Value = '45812 St Mary''s Primary School Hampton'
$ws.Range("B19").Value = 10
$ws.Range("A20").Value = '45967 St Clement of Rome School Bulleen'
$ws.Range("B20").Value = 10
$ws.Range("A21").Value = '46037 Nazareth Catholic Primary School Grovedale'
$ws.Range("B21").Value = 18
$ws.Range("A22").Value = '46050 Our Lady''s Catholic Primary School Craigieburn'
$ws.Range("B22").Value = 36
$ws.Range("A23").Value = '46125 Our Lady of the Southern Cross Primary School Manor Lakes'
$ws.Range("B23").Value = 26
$ws.Range("A24").Value = '46328 Ilim College Dallas Primary Campus Inverloch Cres Tier 1A Dallas'
$ws.Range("B24").Value = 12
$ws.Range("A25").Value = '46390 Al Siraat College Epping'
$ws.Range("B25").Value = 51
$ws.Range("A26").Value = '50681 Broadmeadows Special Developmental School Broadmeadows'
$ws.Range("B26").Value = 10
$ws.Range("A27").Value = '52380 Al Iman College Melton South'
$ws.Range("B27").Value = 22
$ws.Range("A28").Value = '52786 Hume Anglican Grammar Donnybrook Campus'
$ws.Range("B28").Value = 17
$ws.Range("A29").Value = 'Adass Israel School Elsternwick'
$ws.Range("B29").Value = 10
$ws.Range("A30").Value = 'Alfred Health Caulfield Hospital Caulfield'
$ws.Range("B30").Value = 13
$ws.Range("A31").Value = 'Ilim College Dallas Main Campus Dallas Oct'
$ws.Range("B31").Value = 28
$ws.Range("A32").Value = 'Islamic College of Melbourne Tarneit Oct Nov'
$ws.Range("B32").Value = 12
$ws.Range("A33").Value = 'John Henry Primary School Pakenham'
$ws.Range("B33").Value = 10
$ws.Range("A34").Value = 'Melton Willows Melton'
$ws.Range("B34").Value = 10
$ws.Range("A35").Value = 'Morwell Park Primary School Morwell'
$ws.Range("B35").Value = 10
